# Apply updated cryptocurrency market data (prices & 1h volume change)
# to Sheet1 of the workbook, matching the latest scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "42.620.80"
$cell.Style = $origStyle
$ws.Cells.Item(2, 5).Value = "  -7.34%  "
$cell = $ws.Cells.Item(3, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.544.25"
$cell.Style = $origStyle
$ws.Cells.Item(3, 5).Value = "  -1.88%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "295.73"
$cell.Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  -5.15%  "
$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "91.07"
$cell.Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  -7.49%  "
$ws.Cells.Item(7, 5).Value = "  -4.33%  "
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$cell = $ws.Cells.Item(9, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.546"
$cell.Style = $origStyle
$ws.Cells.Item(9, 5).Value = "  -5.65%  "
$cell = $ws.Cells.Item(10, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "35.74"
$cell.Style = $origStyle
$ws.Cells.Item(10, 5).Value = "  -8.04%  "
$ws.Cells.Item(11, 5).Value = "  -4.07%  "
$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.66"
$cell.Style = $origStyle
$ws.Cells.Item(12, 5).Value = "  -5.64%  "
$cell = $ws.Cells.Item(13, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.938.32"
$cell.Style = $origStyle
$ws.Cells.Item(13, 5).Value = "  -1.85%  "
$ws.Cells.Item(14, 5).Value = "  +0.07%  "
$cell = $ws.Cells.Item(15, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.547.13"
$cell.Style = $origStyle
$ws.Cells.Item(15, 5).Value = "  -1.67%  "
$cell = $ws.Cells.Item(16, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.862"
$cell.Style = $origStyle
$ws.Cells.Item(16, 5).Value = "  -5.33%  "
$ws.Cells.Item(17, 5).Value = "  -5.20%  "
$cell = $ws.Cells.Item(18, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "42.618.38"
$cell.Style = $origStyle
$ws.Cells.Item(18, 5).Value = "  -7.66%  "
$cell = $ws.Cells.Item(19, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.61"
$cell.Style = $origStyle
$ws.Cells.Item(19, 5).Value = "  -1.15%  "
$cell = $ws.Cells.Item(20, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0969"
$cell.Style = $origStyle
$ws.Cells.Item(20, 5).Value = "  -4.66%  "
$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.46"
$cell.Style = $origStyle
$ws.Cells.Item(21, 5).Value = "  -2.67%  "
$cell = $ws.Cells.Item(22, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "72.48"
$cell.Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  -0.49%  "
$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "258.73"
$cell.Style = $origStyle
$ws.Cells.Item(23, 5).Value = "  -11.44%  "
$cell = $ws.Cells.Item(24, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.87"
$cell.Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  -6.21%  "
$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "29.37"
$cell.Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  -1.35%  "
$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.10"
$cell.Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  -6.35%  "
$ws.Cells.Item(27, 5).Value = "  +0.06%  "
$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.93"
$cell.Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  -7.29%  "
$ws.Cells.Item(29, 5).Value = "  -4.18%  "
$cell = $ws.Cells.Item(30, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "35.78"
$cell.Style = $origStyle
$ws.Cells.Item(30, 5).Value = "  -5.71%  "
$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.86"
$cell.Style = $origStyle
$ws.Cells.Item(31, 5).Value = "  -5.72%  "
$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "150.06"
$cell.Style = $origStyle
$ws.Cells.Item(32, 5).Value = "  -3.37%  "
$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.15"
$cell.Style = $origStyle
$ws.Cells.Item(33, 5).Value = "  -1.92%  "
$cell = $ws.Cells.Item(34, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.37"
$cell.Style = $origStyle
$ws.Cells.Item(34, 5).Value = "  -5.45%  "
$ws.Cells.Item(35, 5).Value = "  -3.18%  "
$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0789"
$cell.Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  -5.78%  "
$cell = $ws.Cells.Item(37, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$cell.Style = $origStyle
$ws.Cells.Item(37, 5).Value = "  -6.81%  "
$ws.Cells.Item(38, 2).Value = "Stellar"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Cells.Item(38, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.119"
$cell.Style = $origStyle
$ws.Cells.Item(38, 5).Value = "  -3.28%  "
$ws.Cells.Item(39, 2).Value = "EnergySwap"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Cells.Item(39, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "23.92"
$cell.Style = $origStyle
$ws.Cells.Item(39, 5).Value = "  +7.46%  "
$cell = $ws.Cells.Item(40, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.97"
$cell.Style = $origStyle
$ws.Cells.Item(40, 5).Value = "  +1.44%  "
$cell = $ws.Cells.Item(41, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.40"
$cell.Style = $origStyle
$ws.Cells.Item(41, 5).Value = "  -4.65%  "
$cell = $ws.Cells.Item(42, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0307"
$cell.Style = $origStyle
$ws.Cells.Item(42, 5).Value = "  -6.83%  "
$ws.Cells.Item(43, 2).Value = "Maker"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell = $ws.Cells.Item(43, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.072.26"
$cell.Style = $origStyle
$ws.Cells.Item(43, 5).Value = "  -1.06%  "
$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.79"
$cell.Style = $origStyle
$ws.Cells.Item(44, 5).Value = "  -4.06%  "
$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  +0.01%  "
$cell = $ws.Cells.Item(46, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "84.45"
$cell.Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  -13.09%  "
$ws.Cells.Item(47, 5).Value = "  +3.03%  "
$cell = $ws.Cells.Item(48, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.793.49"
$cell.Style = $origStyle
$ws.Cells.Item(48, 5).Value = "  -1.89%  "
$ws.Cells.Item(49, 2).Value = "FraxShare"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Cells.Item(49, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.74"
$cell.Style = $origStyle
$ws.Cells.Item(49, 5).Value = "  -9.61%  "
$ws.Cells.Item(50, 2).Value = "Stacks"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.70"
$cell.Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  -2.79%  "
$cell = $ws.Cells.Item(51, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "102.87"
$cell.Style = $origStyle
$ws.Cells.Item(51, 5).Value = "  -5.03%  "
